$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add columns P and Q with the next sequence values,
# copying the existing bold/centered/bordered style used by the other header cells.
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

$headerStyleSrc = $ws.Range($ws.Cells.Item(1, 15), $ws.Cells.Item(1, 15))
$headerStyleDst = $ws.Range($ws.Cells.Item(1, 16), $ws.Cells.Item(1, 17))
$headerStyleSrc.Copy()
$headerStyleDst.PasteSpecial(-4122)

# --- Rows 2-25: swap the I/K/M/O values and append new P/Q columns.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column
}
